{"js": "// 1. Date change: \"30 de noviembre de 2019\" -> \"01 de diciembre de 2019\"\n//    Replace just the two affected words, leaving the other runs untouched.\nlet results = context.document.body.search(\"30\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"01\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\nresults = context.document.body.search(\"noviembre\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"diciembre\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2. Remove the \"Instituciones y Praxeolog\u00edas\" Heading 3 paragraph entirely.\nresults = context.document.body.search(\"Instituciones y Praxeolog\u00edas\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nfor (const r of results.items) {\n  const paras = r.paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n  paras.items[0].delete();\n  await context.sync();\n}\n\n// 3. Replace the \"As\u00ed, un individuo\" BodyText paragraph with the expanded\n//    explanation, and add a second BodyText paragraph right after it that\n//    folds in (and extends) the old \"La noci\u00f3n de transposici\u00f3n...\" text.\nresults = context.document.body.search(\"As\u00ed, un individuo\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nfor (const r of results.items) {\n  const paras = r.paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n  const p = paras.items[0];\n  p.insertText(\n    \"Los miembros de una instituci\u00f3n son capaces de realizar tipos de tareas que son partes de praxeolog\u00edas puntuales, locales, regionales o globales. Al conjunto de conocimientos conformados de esta manera relacionados con un individuo o una instituci\u00f3n le llamamos equipamiento praxeol\u00f3gico.\",\n    Word.InsertLocation.replace\n  );\n  const newPara = p.insertParagraph(\n    \"La noci\u00f3n de transposici\u00f3n es uno de los m\u00e1s importantes en la teor\u00eda. Esta noci\u00f3n se hace necesaria cuando nos damos cuenta que los equipamientos praxeol\u00f3gicos son modificados cuando transitan de una instituci\u00f3n a otra.\",\n    Word.InsertLocation.after\n  );\n  newPara.style = \"Body Text\";\n  await context.sync();\n}\n\n// 4. Remove the \"Transposici\u00f3n\" Heading 3 paragraph entirely. (The word\n//    \"Transposici\u00f3n\" also appears as a glossary term, so disambiguate by\n//    paragraph style.)\nresults = context.document.body.search(\"Transposici\u00f3n\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\nfor (const r of results.items) {\n  const paras = r.paragraphs;\n  paras.load(\"items/style\");\n  await context.sync();\n  if (paras.items[0].style === \"Heading 3\") {\n    paras.items[0].delete();\n    await context.sync();\n  }\n}\n\n// 5. Remove the old standalone \"La noci\u00f3n de transposici\u00f3n es uno de los\n//    m\u00e1s importantes en la teor\u00eda.\" FirstParagraph paragraph entirely\n//    (its content now lives, expanded, in the new BodyText paragraph from\n//    step 3). Disambiguate by paragraph style since the new paragraph\n//    begins with the same sentence.\nresults = context.document.body.search(\n  \"La noci\u00f3n de transposici\u00f3n es uno de los m\u00e1s importantes en la teor\u00eda.\",\n  { matchCase: true }\n);\nresults.load(\"items\");\nawait context.sync();\nfor (const r of results.items) {\n  const paras = r.paragraphs;\n  paras.load(\"items/style\");\n  await context.sync();\n  if (paras.items[0].style === \"First Paragraph\") {\n    paras.items[0].delete();\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Date change: \"30 de noviembre de 2019\" -> \"01 de diciembre de 2019\"\n#    Replace just the two affected words, leaving the rest of the paragraph intact.\n$r1 = $d.Content\n$r1.Find.ClearFormatting()\n$r1.Find.Text = \"30\"\n$r1.Find.MatchWholeWord = $true\n$r1.Find.MatchCase = $true\n$r1.Find.Forward = $true\n$r1.Find.Wrap = 0\nif ($r1.Find.Execute()) {\n    $r1.Text = \"01\"\n}\n\n$r2 = $d.Content\n$r2.Find.ClearFormatting()\n$r2.Find.Text = \"noviembre\"\n$r2.Find.MatchWholeWord = $true\n$r2.Find.MatchCase = $true\n$r2.Find.Forward = $true\n$r2.Find.Wrap = 0\nif ($r2.Find.Execute()) {\n    $r2.Text = \"diciembre\"\n}\n\n# 2. Remove the \"Instituciones y Praxeolog\u00edas\" Heading 3 paragraph entirely.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $txt = $p.Range.Text.Trim()\n    $style = $p.Range.Style.NameLocal\n    if ($style -eq \"Heading 3\" -and $txt -eq \"Instituciones y Praxeolog\u00edas\") {\n        $p.Range.Delete()\n    }\n}\n\n# 3. Replace the \"As\u00ed, un individuo\" BodyText paragraph with the expanded\n#    explanation, followed immediately by a second BodyText paragraph that\n#    folds in (and extends) the old \"La noci\u00f3n de transposici\u00f3n...\" text.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $txt = $p.Range.Text.Trim()\n    $style = $p.Range.Style.NameLocal\n    if ($style -eq \"Body Text\" -and $txt -eq \"As\u00ed, un individuo\") {\n        $p.Range.Text = \"Los miembros de una instituci\u00f3n son capaces de realizar tipos de tareas que son partes de praxeolog\u00edas puntuales, locales, regionales o globales. Al conjunto de conocimientos conformados de esta manera relacionados con un individuo o una instituci\u00f3n le llamamos equipamiento praxeol\u00f3gico.`rLa noci\u00f3n de transposici\u00f3n es uno de los m\u00e1s importantes en la teor\u00eda. Esta noci\u00f3n se hace necesaria cuando nos damos cuenta que los equipamientos praxeol\u00f3gicos son modificados cuando transitan de una instituci\u00f3n a otra.\"\n        break\n    }\n}\n\n# 4. Remove the \"Transposici\u00f3n\" Heading 3 paragraph, and the old standalone\n#    \"La noci\u00f3n de transposici\u00f3n es uno de los m\u00e1s importantes en la\n#    teor\u00eda.\" FirstParagraph paragraph, entirely (its content now lives,\n#    expanded, in the new BodyText paragraph created in step 3). Walk\n#    backwards so deleting doesn't shift the indices we still need to\n#    visit. Style disambiguates these from the similarly-worded glossary\n#    entry and the new BodyText paragraph.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $txt = $p.Range.Text.Trim()\n    $style = $p.Range.Style.NameLocal\n    if ($style -eq \"Heading 3\" -and $txt -eq \"Transposici\u00f3n\") {\n        $p.Range.Delete()\n    }\n    elseif ($style -eq \"First Paragraph\" -and $txt -eq \"La noci\u00f3n de transposici\u00f3n es uno de los m\u00e1s importantes en la teor\u00eda.\") {\n        $p.Range.Delete()\n    }\n}\n"}
